$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: the phone number cell was re-typed as a genuine number (was text)
$ws.Cells.Item(35, 1).Value = 71277620

# Row 36: new payment record 71277620 (Cash) 2025-08-18T17:10:26
# Phone number looks numeric but must stay text -> force with a leading
# apostrophe, then strip the resulting quote-prefix formatting back to Normal.
$ws.Cells.Item(36, 1).Value = "'71277620"
$ws.Cells.Item(36, 1).Style = "Normal"

# Blank text cells (phone-less placeholder columns) - force text type via
# a bare apostrophe so they round-trip as empty strings, not blank cells.
$ws.Cells.Item(36, 2).Value = "'"
$ws.Cells.Item(36, 2).Style = "Normal"

$ws.Cells.Item(36, 3).Value = "Cash"
$ws.Cells.Item(36, 4).Value = "2025-08-18T17:10:26"
$ws.Cells.Item(36, 5).Value = 76

$ws.Cells.Item(36, 6).Value = "'"
$ws.Cells.Item(36, 6).Style = "Normal"

$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 76
